$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H3").Value2 = 6.2
$ws.Range("I3").Value2 = 7.8
$ws.Range("P3").Value2 = 1.99
$ws.Range("S3").Value2 = 3.1
$ws.Range("X3").Value2 = 1000
$ws.Range("AC3").Value2 = 970
$ws.Range("AG3").Value2 = 970
$ws.Range("AK3").Value2 = 970
$ws.Range("AL3").Value2 = 1000
$ws.Range("AN3").Value2 = 970
$ws.Range("G4").Value2 = 5.9
$ws.Range("M4").Value2 = 1.05
$ws.Range("Q4").Value2 = 1.68
$ws.Range("R4").Value2 = 1.56
$ws.Range("S4").Value2 = 2.7
$ws.Range("W4").Value2 = 1.2
$ws.Range("AA4").Value2 = 16
$ws.Range("AD6").Value2 = 23
$ws.Range("AF6").Value2 = 16
$ws.Range("AG6").Value2 = 12.5
$ws.Range("AN6").Value2 = 7.4
$ws.Range("F7").Value2 = 8.2
$ws.Range("G7").Value2 = 8.6
$ws.Range("K7").Value2 = 6.6
$ws.Range("Y7").Value2 = 15.5
$ws.Range("Z8").Value2 = 970
$ws.Range("AH8").Value2 = 970
$ws.Range("AL8").Value2 = 55
$ws.Range("K10").Value2 = 950
$ws.Range("N10").Value2 = 3.55
$ws.Range("P10").Value2 = 1.86
$ws.Range("Q10").Value2 = 1.8
$ws.Range("S10").Value2 = 2.78
$ws.Range("X10").Value2 = 970
$ws.Range("Z10").Value2 = 970
$ws.Range("AA10").Value2 = 970
$ws.Range("AC10").Value2 = 970
$ws.Range("AD10").Value2 = 970
$ws.Range("AF10").Value2 = 70
$ws.Range("AO10").Value2 = 970
$ws.Range("F11").Value2 = 1.55
$ws.Range("G11").Value2 = 1.78
$ws.Range("I11").Value2 = 7.4
$ws.Range("J11").Value2 = 3.4
$ws.Range("L11").Value2 = 1.19
$ws.Range("N11").Value2 = 5.4
$ws.Range("P11").Value2 = 2.5
$ws.Range("R11").Value2 = 1.61
$ws.Range("T11").Value2 = 1.49
$ws.Range("W11").Value2 = 2.28
$ws.Range("X11").Value2 = 1000
$ws.Range("Y11").Value2 = 1000
$ws.Range("Z11").Value2 = 1000
$ws.Range("AB11").Value2 = 1000
$ws.Range("AC11").Value2 = 1000
$ws.Range("AD11").Value2 = 1000
$ws.Range("AF11").Value2 = 1000
$ws.Range("AG11").Value2 = 1000
$ws.Range("AH11").Value2 = 1000
$ws.Range("AJ11").Value2 = 1000
$ws.Range("AK11").Value2 = 1000
$ws.Range("AL11").Value2 = 1000
$ws.Range("AN11").Value2 = 1000
$ws.Range("AO11").Value2 = 1000
$ws.Range("J13").Value2 = 3.7
$ws.Range("P13").Value2 = 2.06
$ws.Range("Q13").Value2 = 1.74
$ws.Range("F14").Value2 = 1.35
$ws.Range("G14").Value2 = 1.58
$ws.Range("H14").Value2 = 6.8
$ws.Range("I14").Value2 = 12.5
$ws.Range("J14").Value2 = 4.4
$ws.Range("R14").Value2 = 1.42
$ws.Range("S14").Value2 = 2.44
$ws.Range("T14").Value2 = 1.71
$ws.Range("U14").Value2 = 1.65
$ws.Range("V14").Value2 = 1.09
$ws.Range("W14").Value2 = 2.72
$ws.Range("X14").Value2 = 28
$ws.Range("Y14").Value2 = 40
$ws.Range("Z14").Value2 = 100
$ws.Range("AB14").Value2 = 12
$ws.Range("AD14").Value2 = 46
$ws.Range("AH14").Value2 = 36
$ws.Range("F17").Value2 = 4.6
$ws.Range("G17").Value2 = 4.8
$ws.Range("H17").Value2 = 1.92
$ws.Range("I17").Value2 = 1.93
$ws.Range("P17").Value2 = 1.86
$ws.Range("Q17").Value2 = 2.14
$ws.Range("S17").Value2 = 3.9
$ws.Range("V17").Value2 = 2.06
$ws.Range("AO17").Value2 = 15
$ws.Range("T19").Value2 = 1.64
$ws.Range("I20").Value2 = 5.8
$ws.Range("K20").Value2 = 5
$ws.Range("U20").Value2 = 2.62
$ws.Range("AI20").Value2 = 50
$ws.Range("AO20").Value2 = 40
$ws.Range("F21").Value2 = 1.29
$ws.Range("G21").Value2 = 1.3
$ws.Range("I21").Value2 = 12
$ws.Range("J21").Value2 = 6.8
$ws.Range("R21").Value2 = 1.98
$ws.Range("T21").Value2 = 1.75
$ws.Range("W21").Value2 = 4.3
$ws.Range("Z21").Value2 = 130
$ws.Range("AA21").Value2 = 1000
$ws.Range("AI21").Value2 = 1000
$ws.Range("AJ21").Value2 = 11.5
$ws.Range("AM21").Value2 = 1000